$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension implicitly handled by Excel; set data rows 2-7

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Bdnf"
$ws.Range("C2").Value = "Ngfr"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.020961333333334
$ws.Range("H2").Value = 6.062884
$ws.Range("I2").Value = 0.3447258214530571
$ws.Range("J2").Value = 0.3447258214530571
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.6946430000000001
$ws.Range("N2").Value = 2.083929
$ws.Range("O2").Value = 0.1269399741689062
$ws.Range("P2").Value = 0.1269399741689062
$ws.Range("Q2").Value = 1.403846643470667
$ws.Range("R2").Value = 12.634619791236
$ws.Range("S2").Value = 0.04375948687060605
$ws.Range("T2").Value = 0.04375948687060605

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Bdnf"
$ws.Range("C3").Value = "Ngfr"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.020961333333334
$ws.Range("H3").Value = 6.062884
$ws.Range("I3").Value = 0.3447258214530571
$ws.Range("J3").Value = 0.3447258214530571
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05042666666666667
$ws.Range("N3").Value = 0.15128
$ws.Range("O3").Value = 0.00921503529739839
$ws.Range("P3").Value = 0.00921503529739839
$ws.Range("Q3").Value = 0.1019103435022222
$ws.Range("R3").Value = 0.9171930915200001
$ws.Range("S3").Value = 0.003176660612614576
$ws.Range("T3").Value = 0.003176660612614576

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Bdnf"
$ws.Range("C4").Value = "Ngfr"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.020961333333334
$ws.Range("H4").Value = 6.062884
$ws.Range("I4").Value = 0.3447258214530571
$ws.Range("J4").Value = 0.3447258214530571
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.727146666666667
$ws.Range("N4").Value = 14.18144
$ws.Range("O4").Value = 0.8638449905336953
$ws.Range("P4").Value = 0.8638449905336955
$ws.Range("Q4").Value = 9.55338063032889
$ws.Range("R4").Value = 85.98042567296001
$ws.Range("S4").Value = 0.2977896739698365
$ws.Range("T4").Value = 0.2977896739698365

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Bdnf"
$ws.Range("C5").Value = "Ngfr"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.841556666666667
$ws.Range("H5").Value = 11.52467
$ws.Range("I5").Value = 0.6552741785469429
$ws.Range("J5").Value = 0.6552741785469429
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6946430000000001
$ws.Range("N5").Value = 2.083929
$ws.Range("O5").Value = 0.1269399741689062
$ws.Range("P5").Value = 0.1269399741689062
$ws.Range("Q5").Value = 2.668510447603334
$ws.Range("R5").Value = 24.01659402843001
$ws.Range("S5").Value = 0.08318048729830019
$ws.Range("T5").Value = 0.08318048729830019

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Bdnf"
$ws.Range("C6").Value = "Ngfr"
$ws.Range("D6").Value = "M2"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.841556666666667
$ws.Range("H6").Value = 11.52467
$ws.Range("I6").Value = 0.6552741785469429
$ws.Range("J6").Value = 0.6552741785469429
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.05042666666666667
$ws.Range("N6").Value = 0.15128
$ws.Range("O6").Value = 0.00921503529739839
$ws.Range("P6").Value = 0.00921503529739839
$ws.Range("Q6").Value = 0.1937168975111111
$ws.Range("R6").Value = 1.7434520776
$ws.Range("S6").Value = 0.006038374684783814
$ws.Range("T6").Value = 0.006038374684783814

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Bdnf"
$ws.Range("C7").Value = "Ngfr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.841556666666667
$ws.Range("H7").Value = 11.52467
$ws.Range("I7").Value = 0.6552741785469429
$ws.Range("J7").Value = 0.6552741785469429
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.727146666666667
$ws.Range("N7").Value = 14.18144
$ws.Range("O7").Value = 0.8638449905336953
$ws.Range("P7").Value = 0.8638449905336955
$ws.Range("Q7").Value = 18.15960179164445
$ws.Range("R7").Value = 163.4364161248
$ws.Range("S7").Value = 0.5660553165638589
$ws.Range("T7").Value = 0.566055316563859
